$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item('Sheet 1')
$wsData.Range('E5').Value = 'Based on official disease reports to the WOAH'
$wsData.Range('E6').Value = 'Sheep Pox and Goat pox (SGPox) are viral diseases listed in the World Organisation for Animal Health ({ref009:WOAH}) Terrestrial Animal Health Code and must be reported to the WOAH. The map to the right displays outbreak points reported to the WOAH early warning system since 2005.'
$wsData.Range('E7').Value = 'As described in the WOAH {ref005:Terrestrial Animal Health Code}, the WOAH early warning system includes immediate notifications and follow-up reports on:'
$wsData.Range('E14').Value = 'Countries are coloured according to the available information regarding their stable disease situation (disease status legend). This information is provided by countries through the WOAH monitoring system, which is a different reporting channel.<br>Immediate notifications (points) and disease status (country/region colours) are reported to the WOAH in different spatial and temporal scales, and therefore are displayed in the map as layers which can be filtered independently.'
$wsData.Range('E17').Value = 'For more up to date reports, visit the original data source: {ref001:WOAH-WAHIS}.'
$wsData.Range('E21').Value = 'A summary of the disease in animal hosts is given in the {ref008:WOAH Technical disease card}.'
$wsData.Range('E34').Value = 'Humans are not susceptible to SGPoxV and therefore there is no direct impact on public health {ref008:WOAH Technical disease card} .'
$wsData.Range('E43').Value = 'Refer to the {ref008:WOAH Technical disease card} for a key summary of the virus characteristics. '
$wsData.Range('E55').Value = 'Refer to the {ref008:WOAH Technical disease card} for a key summary of the disease transmission and epidemiological parameters.'
$wsData.Range('E68').Value = 'WOAH-prescribed tests for international trade include:the commercial double-antigen enzyme-linked immunosorbent assay (ELISA) and (real-time) polymerase chain reaction (PCR) excluding vaccine strains ({ref010:WOAH, Terrestriam Manual},{ref034:Haegeman et al. 2020})'
$wsData.Range('E94').Value = 'Geographical distribution data has been kindly provided by the World Organisation of Animal Health (WOAH). {ref001:WOAH-WAHIS} (WOAH World Animal Health Information System) is the original source of these data.'
$wsRefs = $wb.Worksheets.Item('References')
$wsRefs.Range('C2').Value = 'WOAH-WAHIS (WOAH World Animal Health Information System)'
$wsRefs.Range('C5').Value = 'WOAH (World Organisation for Animal Health). Terrestrial Animal Health Code 2021. WOAH, Paris, France'
$wsRefs.Range('C8').Value = 'WOAH (World Organisation for Animal Health) Technical Disease Card: Sheep pox and goat pox. 2013.'
$wsRefs.Range('C9').Value = 'WOAH (World Organisation for Animal Health), 2018.Sheep pox and goat pox. Chapter 3.7.12. WOAH Terrestrial Animal Health Code 2018. WOAH, Paris, France'
$wsRefs.Range('C10').Value = 'WOAH (World Organisation for Animal Health), 2017. Sheep and Goat Pox. Chapter 2.07.13. WOAH Terrestrial Manual. WOAH, Paris, France'
